# "very reduced model runs now"
#
# On the "Constants" sheet, the "Treatment recovery on unsuppressive ART"
# block (which only held a single "People on unsuppressive ART who recover (%)"
# row) is replaced by a "CD4 change due to non-suppressive ART (% per year)"
# block that lists all 8 CD4-stage transitions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")
$ws.Activate()

# Make room: the new block needs 8 data rows (45-52) instead of the single
# old row 45, so insert 7 new rows right after the existing row 45.
$ws.Range("A46:E52").EntireRow.Insert()

# Section heading
$ws.Range("A43").Value = "CD4 change due to non-suppressive ART (% per year)"

# New data rows: label (column B) + best/low/high (columns C/D/E)
$labels = @(
  "CD4(500) to CD4(350-500)",
  "CD4(350-500) to CD4(>500)",
  "CD4(350-500) to CD4(200-350)",
  "CD4(200-350) to CD4(350-500)",
  "CD4(200-350) to CD4(50-200)",
  "CD4(50-200) to CD4(200-350)",
  "CD4(50-200) to CD4(<50)",
  "CD4(<50) to CD4(50-200)"
)
$best = @(0.026, 0.15, 0.1, 0.053, 0.162, 0.117, 0.09, 0.111)
$low  = @(0.005, 0.038, 0.022, 0.008, 0.05, 0.032, 0.019, 0.047)
$high = @(0.275, 0.885, 0.87, 0.827, 0.869, 0.686, 0.723, 0.563)

for ($i = 0; $i -lt $labels.Length; $i++) {
  $r = 45 + $i
  $ws.Range("B$r").Value = $labels[$i]
  $ws.Range("C$r").Value = $best[$i]
  $ws.Range("D$r").Value = $low[$i]
  $ws.Range("E$r").Value = $high[$i]
}

$ws.Range("C45:E52").NumberFormat = "0.0%"

$ws.Range("B55").Select()

Write-Host "done"
